$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 1 (header labels): GE1:HP1 gets relabeled/shifted values ---
$row1Values = @(
    't_EOc0',
    't_EOc1',
    't_EOc2',
    't_EOc3',
    't_to10',
    't_to11',
    't_to12',
    't_to13',
    't_to20',
    't_to21',
    't_to22',
    't_to23',
    't_tr10000',
    't_tr10001',
    't_tr10002',
    't_tr10003',
    'T_EO_year0',
    'T_EO_year1',
    'T_EO_year2',
    'T_EO_year3',
    'T_to1_year0',
    'T_to1_year1',
    'T_to1_year2',
    'T_to1_year3',
    'T_to2_year0',
    'T_to2_year1',
    'T_to2_year2',
    'T_to2_year3',
    'T_tr1000_year0',
    'T_tr1000_year1',
    'T_tr1000_year2',
    'T_tr1000_year3',
    'T_EO_year_all',
    'T_D1_year_all',
    'T_D2_year_all',
    'T_to1_year_all',
    'T_to2_year_all',
    'T_tr1000_year_all'
)

# --- Row 2 (data): GE2:HP2 gets recomputed numeric values ---
$row2Values = @(
    0.31,
    0.3,
    0.44,
    0.2,
    5.74,
    3.96,
    9.08,
    5.24,
    22.78,
    15.84,
    29.04,
    20.94,
    4.59,
    4.28,
    9.58,
    2.23,
    2123.17,
    9334.969999999999,
    15053.7,
    2970.09,
    922.1900000000001,
    2887.32,
    8315.83,
    2918.31,
    1109.84,
    3614.05,
    8238.07,
    3640.42,
    3726.95,
    11562.9,
    35381.79,
    3802.11,
    29481.93,
    1920.11,
    1935.94,
    13668.28,
    15211.18,
    53384.27
)

$colCount = $row1Values.Count

$arr1 = New-Object "object[,]" 1,$colCount
$arr2 = New-Object "object[,]" 1,$colCount
for ($i = 0; $i -lt $colCount; $i++) {
    $arr1[0,$i] = $row1Values[$i]
    $arr2[0,$i] = [double]$row2Values[$i]
}

$ws.Range("GE1:HP1").Value = $arr1
$ws.Range("GE2:HP2").Value = $arr2

# --- Remove the now-obsolete trailing columns (old HQ:IK) entirely ---
$ws.Range("HQ1:IK2").ClearContents()

